# FormatoRecibo.docx - "cambios en flujo de ingresos"
#
# The address line in the receipt header has a misspelled district name:
# "Occlo" must read "Ocllo" ("Urb. Ocllo" is the real neighbourhood name).
# Replace the whole "Av. Chimpu Occlo N°198 Urb. " run of text with the
# corrected spelling; the trailing "Lucyana" word (and its spell-check
# markers) is left untouched, exactly like in the source edit.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Av. Chimpu Occlo N" + [char]0x00B0 + "198 Urb. ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Av. Chimpu Ocllo N" + [char]0x00B0 + "198 Urb. ",
    2
)
